# "Dodano trudne pytanie i napisano zdjecie"
# Adds a new "hard" (poz_trud = 3) question about SDG logos to Arkusz1,
# renames several existing "kat" (category) values, widens column D,
# and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Rename "kat" values for the existing answer groups -------------------
# Category 1 (rows 2-5): ge -> npov
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 2).Value = "npov"
}

# Category 3 (rows 6-9): ge -> dwaeg
for ($r = 6; $r -le 9; $r++) {
    $ws.Cells.Item($r, 2).Value = "dwaeg"
}

# Category 4 (rows 10-13): ge -> gen
for ($r = 10; $r -le 13; $r++) {
    $ws.Cells.Item($r, 2).Value = "gen"
}

# Category 5 (rows 14-17) is left untouched (stays "ge").

# Category 2 (rows 18-23): ge -> gen
for ($r = 18; $r -le 23; $r++) {
    $ws.Cells.Item($r, 2).Value = "gen"
}

# --- Append the new "hard" question (rows 24-31) ---------------------------
$answers = @(
    "Affordable and clean energy",
    "Life on land ",
    "Zero hunger",
    "No poverty",
    "Climate action",
    "Responsible consumption and production",
    "Life below water",
    "Quality education"
)
$correct = @("F", "T", "F", "F", "F", "F", "F", "F")

for ($i = 0; $i -lt $answers.Length; $i++) {
    $r = 24 + $i
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "gen"
    $ws.Cells.Item($r, 3).Value = 3
    $ws.Cells.Item($r, 4).Value = "Which SDG Logo you saw? "
    $ws.Cells.Item($r, 5).Value = $answers[$i]
    $ws.Cells.Item($r, 6).Value = $correct[$i]
    $ws.Cells.Item($r, 7).Value = "F"
    $ws.Cells.Item($r, 8).Value = "pyt_trud.png"
}

# --- Formatting / view tweaks ----------------------------------------------
# Widen column D (pytanie) slightly to fit the new, longer question text.
$ws.Columns.Item(4).ColumnWidth = 52.5

# Move the selection to reflect where editing finished (column E, new last row).
$ws.Range("E32").Select() | Out-Null
